$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran the averaged-intensities notebook after adding three new spiral
# sampling schemes ("Spiral-90deg-10rot-5space", "Spiral-90deg-15rot-5space",
# "Spiral-90deg-10rot-3space") ahead of the hex-grid schemes in the list of
# orientation schemes, and moving "Gaussian-Quadrature" earlier in that same
# list. The per-scheme results table on the sheet grows by three rows (to
# r19) and the scheme labels for the existing rows 10-16 shift down to match
# the new scheme ordering.

# Make the three brand new rows look like the existing index rows: copy the
# bold/bordered "HKL index" style from A16 down onto A17:A19 first.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

# Rows 3-9 (ND Single ... Ring Perpendicular to TD) are unaffected by the
# re-ordering, so only rows 10 onward need their scheme label (column B)
# rewritten; column A (the HKL row index) stays sequential 0..17.

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10:M10").Value = 1

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11:M11").Value = 1

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12:M12").Value = 1

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13:M13").Value = 1

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14:M14").Value = 1

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15:M15").Value = 1

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16:M16").Value = 1

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17:M17").Value = 1

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18:M18").Value = 1

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19:M19").Value = 1
